$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row: "_old" columns -> "_FV2210", "_new" columns -> "_FV2304" ---
# Columns A:J carry the "*_old" headers, column K is "diff" (unchanged),
# columns L:U carry the "*_new" headers.
$leftHeaders  = @("Segmentname_FV2210","Segmentgruppe_FV2210","Segment_FV2210","Datenelement_FV2210","Segment ID_FV2210","Code_FV2210","Qualifier_FV2210","Beschreibung_FV2210","Bedingungsausdruck_FV2210","Bedingung_FV2210")
$rightHeaders = @("Segmentname_FV2304","Segmentgruppe_FV2304","Segment_FV2304","Datenelement_FV2304","Segment ID_FV2304","Code_FV2304","Qualifier_FV2304","Beschreibung_FV2304","Bedingungsausdruck_FV2304","Bedingung_FV2304")

for ($i = 0; $i -lt $leftHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $leftHeaders[$i]
}
for ($i = 0; $i -lt $rightHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $rightHeaders[$i]
}

# --- 2. Turn the used range into an Excel Table ("Table1") with an AutoFilter ---
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:U56"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (split below row 1, frozen) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
